# Auto-generated edit script applying the Garuda_Profits diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 30505.666
$ws.Cells.Item(21, 9).Value = 25758.5
$ws.Cells.Item(21, 11).Value = 25758.5
$ws.Cells.Item(21, 13).Value = -25290.5

$ws.Cells.Item(23, 8).Value = 30505.666
$ws.Cells.Item(23, 9).Value = 25758.5
$ws.Cells.Item(23, 11).Value = 25758.5
$ws.Cells.Item(23, 13).Value = -25524.5

$ws.Cells.Item(40, 8).Value = 2120.389
$ws.Cells.Item(40, 9).Value = 3994.5
$ws.Cells.Item(40, 10).Value = 1584.9286
$ws.Cells.Item(40, 11).Value = 3994.5
$ws.Cells.Item(40, 12).Value = 1584.9286
$ws.Cells.Item(40, 13).Value = -3819.5
$ws.Cells.Item(40, 14).Value = -1934.9286

$ws.Cells.Item(98, 8).Value = 8833.27
$ws.Cells.Item(98, 9).Value = 9147
$ws.Cells.Item(98, 11).Value = 9147
$ws.Cells.Item(98, 13).Value = -7649

$ws.Cells.Item(122, 8).Value = 8833.27
$ws.Cells.Item(122, 9).Value = 9147
$ws.Cells.Item(122, 11).Value = 27441
$ws.Cells.Item(122, 13).Value = -24991

$ws.Cells.Item(129, 8).Value = 1858540
$ws.Cells.Item(129, 9).Value = 556.75
$ws.Cells.Item(129, 10).Value = 3344926.5
$ws.Cells.Item(129, 11).Value = 1670.25
$ws.Cells.Item(129, 12).Value = 10034779.5
$ws.Cells.Item(129, 13).Value = 3329.75
$ws.Cells.Item(129, 14).Value = -10044779.5

$ws.Cells.Item(133, 8).Value = 41100
$ws.Cells.Item(133, 10).Value = 48500
$ws.Cells.Item(133, 12).Value = 48500
$ws.Cells.Item(133, 14).Value = -58620

$ws.Cells.Item(134, 8).Value = 37228.57
$ws.Cells.Item(134, 10).Value = 37228.57
$ws.Cells.Item(134, 12).Value = 37228.57
$ws.Cells.Item(134, 14).Value = -47368.57

$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 14).ClearContents()

$ws.Cells.Item(137, 8).Value = 33335322
$ws.Cells.Item(137, 9).Value = 1836.7273
$ws.Cells.Item(137, 11).Value = 5510.1819
$ws.Cells.Item(137, 13).Value = -2960.1819

$ws.Cells.Item(138, 8).Value = 2624.8774
$ws.Cells.Item(138, 9).Value = 1197.091
$ws.Cells.Item(138, 10).Value = 3349.754
$ws.Cells.Item(138, 11).Value = 3591.273
$ws.Cells.Item(138, 12).Value = 10049.262
$ws.Cells.Item(138, 13).Value = 1548.727
$ws.Cells.Item(138, 14).Value = -20329.262

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 32950.79
$ws.Cells.Item(32, 9).Value = 31593.732
$ws.Cells.Item(32, 11).Value = 31593.732
$ws.Cells.Item(32, 13).Value = -31306.732

$ws.Cells.Item(61, 8).Value = 2339.0833
$ws.Cells.Item(61, 9).Value = 1928.3158
$ws.Cells.Item(61, 10).Value = 3900
$ws.Cells.Item(61, 11).Value = 1928.3158
$ws.Cells.Item(61, 12).Value = 3900
$ws.Cells.Item(61, 13).Value = -1716.3158
$ws.Cells.Item(61, 14).Value = -4324

$ws.Cells.Item(63, 8).Value = 1433405.4
$ws.Cells.Item(63, 9).Value = 2502459.5
$ws.Cells.Item(63, 10).Value = 8000
$ws.Cells.Item(63, 11).Value = 2502459.5
$ws.Cells.Item(63, 12).Value = 8000
$ws.Cells.Item(63, 13).Value = -2501773.5
$ws.Cells.Item(63, 14).Value = -9372

$ws.Cells.Item(66, 8).Value = 1433405.4
$ws.Cells.Item(66, 9).Value = 2502459.5
$ws.Cells.Item(66, 10).Value = 8000
$ws.Cells.Item(66, 11).Value = 12512297.5
$ws.Cells.Item(66, 12).Value = 40000
$ws.Cells.Item(66, 13).Value = -12508865.5
$ws.Cells.Item(66, 14).Value = -46864

$ws.Cells.Item(136, 8).Value = 2339.0833
$ws.Cells.Item(136, 9).Value = 1928.3158
$ws.Cells.Item(136, 10).Value = 3900
$ws.Cells.Item(136, 11).Value = 5784.9474
$ws.Cells.Item(136, 12).Value = 11700
$ws.Cells.Item(136, 13).Value = -3234.9474
$ws.Cells.Item(136, 14).Value = -16800

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 737.5
$ws.Cells.Item(99, 9).Value = 550
$ws.Cells.Item(99, 10).Value = 925
$ws.Cells.Item(99, 11).Value = 550
$ws.Cells.Item(99, 12).Value = 925
$ws.Cells.Item(99, 13).Value = 948
$ws.Cells.Item(99, 14).Value = -3921

$ws.Cells.Item(134, 8).Value = 4346.305
$ws.Cells.Item(134, 9).Value = 4685.095
$ws.Cells.Item(134, 10).Value = 3509.2942
$ws.Cells.Item(134, 11).Value = 14055.285
$ws.Cells.Item(134, 12).Value = 10527.8826
$ws.Cells.Item(134, 13).Value = -11520.285
$ws.Cells.Item(134, 14).Value = -15597.8826

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 19612986
$ws.Cells.Item(31, 9).Value = 2578.9167
$ws.Cells.Item(31, 11).Value = 2578.9167
$ws.Cells.Item(31, 13).Value = -2283.9167

$ws.Cells.Item(34, 8).Value = 19612986
$ws.Cells.Item(34, 9).Value = 2578.9167
$ws.Cells.Item(34, 11).Value = 2578.9167
$ws.Cells.Item(34, 13).Value = -2376.9167

$ws.Cells.Item(52, 8).Value = 32390
$ws.Cells.Item(52, 10).Value = 32390
$ws.Cells.Item(52, 12).Value = 32390
$ws.Cells.Item(52, 14).Value = -32978

$ws.Cells.Item(99, 8).Value = 2014.9429
$ws.Cells.Item(99, 9).Value = 1638.65
$ws.Cells.Item(99, 11).Value = 1638.65
$ws.Cells.Item(99, 13).Value = -140.6500000000001

$ws.Cells.Item(126, 8).Value = 2014.9429
$ws.Cells.Item(126, 9).Value = 1638.65
$ws.Cells.Item(126, 11).Value = 4915.950000000001
$ws.Cells.Item(126, 13).Value = -2445.950000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 680.7292
$ws.Cells.Item(113, 9).Value = 592.1579
$ws.Cells.Item(113, 10).Value = 1017.3
$ws.Cells.Item(113, 11).Value = 1776.4737
$ws.Cells.Item(113, 12).Value = 3051.9
$ws.Cells.Item(113, 13).Value = 393.5263
$ws.Cells.Item(113, 14).Value = -7391.9

$ws.Cells.Item(129, 8).Value = 12821407
$ws.Cells.Item(129, 9).Value = 513.1111
$ws.Cells.Item(129, 11).Value = 1539.3333
$ws.Cells.Item(129, 13).Value = 3460.6667

$ws.Cells.Item(131, 8).Value = 720.11
$ws.Cells.Item(131, 10).Value = 795.4235
$ws.Cells.Item(131, 12).Value = 2386.2705
$ws.Cells.Item(131, 14).Value = -12466.2705

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2807.5667
$ws.Cells.Item(126, 9).Value = 3161.35
$ws.Cells.Item(126, 10).Value = 2100
$ws.Cells.Item(126, 11).Value = 9484.049999999999
$ws.Cells.Item(126, 12).Value = 6300
$ws.Cells.Item(126, 13).Value = -7014.049999999999
$ws.Cells.Item(126, 14).Value = -11240

$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1462.9
$ws.Cells.Item(46, 9).Value = 800
$ws.Cells.Item(46, 10).Value = 1579.8823
$ws.Cells.Item(46, 11).Value = 800
$ws.Cells.Item(46, 12).Value = 1579.8823
$ws.Cells.Item(46, 13).Value = -612
$ws.Cells.Item(46, 14).Value = -1955.8823

$ws.Cells.Item(93, 8).Value = 2075.2856
$ws.Cells.Item(93, 9).Value = 1937.8334
$ws.Cells.Item(93, 10).Value = 2900
$ws.Cells.Item(93, 11).Value = 1937.8334
$ws.Cells.Item(93, 12).Value = 2900
$ws.Cells.Item(93, 13).Value = -689.8334
$ws.Cells.Item(93, 14).Value = -5396

$ws.Cells.Item(122, 8).Value = 4629.7144
$ws.Cells.Item(122, 9).Value = 6802.6665
$ws.Cells.Item(122, 10).Value = 3000
$ws.Cells.Item(122, 11).Value = 20407.9995
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = -17957.9995
$ws.Cells.Item(122, 14).Value = -13900

$ws.Cells.Item(132, 8).Value = 9505.321
$ws.Cells.Item(132, 9).Value = 10658.739
$ws.Cells.Item(132, 11).Value = 31976.217
$ws.Cells.Item(132, 13).Value = -29446.217

$ws.Cells.Item(136, 8).Value = 13697.8
$ws.Cells.Item(136, 9).Value = 18282.572
$ws.Cells.Item(136, 10).Value = 3000
$ws.Cells.Item(136, 11).Value = 54847.716
$ws.Cells.Item(136, 12).Value = 9000
$ws.Cells.Item(136, 13).Value = -52297.716
$ws.Cells.Item(136, 14).Value = -14100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 3685.1
$ws.Cells.Item(122, 9).Value = 3943.0715
$ws.Cells.Item(122, 10).Value = 3083.1667
$ws.Cells.Item(122, 11).Value = 11829.2145
$ws.Cells.Item(122, 12).Value = 9249.500100000001
$ws.Cells.Item(122, 13).Value = -9379.2145
$ws.Cells.Item(122, 14).Value = -14149.5001

$ws.Cells.Item(123, 8).Value = 23047
$ws.Cells.Item(123, 10).Value = 23047
$ws.Cells.Item(123, 12).Value = 23047
$ws.Cells.Item(123, 14).Value = -32847

$ws.Cells.Item(125, 8).Value = 36021.668
$ws.Cells.Item(125, 10).Value = 36021.668
$ws.Cells.Item(125, 12).Value = 36021.668
$ws.Cells.Item(125, 14).Value = -45861.668

$ws.Cells.Item(136, 8).Value = 1734
$ws.Cells.Item(136, 9).Value = 1670
$ws.Cells.Item(136, 10).Value = 2950
$ws.Cells.Item(136, 11).Value = 5010
$ws.Cells.Item(136, 12).Value = 8850
$ws.Cells.Item(136, 13).Value = -2460
$ws.Cells.Item(136, 14).Value = -13950
